{"js": "// Commit: \"Added nuget restor & inline code\"\n//\n// Adds a new character style \"InlineCodeSnippet\" \u2014 based on the existing\n// \"CodeSnippetZchn\" (the linked character style of the paragraph style\n// \"CodeSnippet\") \u2014 so inline code spans can be formatted consistently with\n// the block CodeSnippet style. Mirrors:\n//   <w:style w:type=\"character\" w:customStyle=\"1\" w:styleId=\"InlineCodeSnippet\">\n//     <w:name w:val=\"InlineCodeSnippet\"/>\n//     <w:basedOn w:val=\"CodeSnippetZchn\"/>\n//     <w:uiPriority w:val=\"1\"/>\n//     <w:qFormat/>\n//   </w:style>\n\n// Create the new character style.\ncontext.document.addStyle(\"InlineCodeSnippet\", Word.StyleType.character);\nawait context.sync();\n\n// Re-fetch it by name so the follow-up property writes land on a fully\n// resolved style object (rather than the anchor returned inline from\n// addStyle, before it has been synced back).\nconst styles = context.document.getStyles();\nconst inlineCodeSnippet = styles.getByName(\"InlineCodeSnippet\");\n\n// basedOn -> <w:basedOn w:val=\"CodeSnippetZchn\"/>\ninlineCodeSnippet.baseStyle = \"CodeSnippetZchn\";\n// uiPriority -> <w:uiPriority w:val=\"1\"/>\ninlineCodeSnippet.priority = 1;\n// qFormat -> <w:qFormat/>\ninlineCodeSnippet.quickStyle = true;\n\nawait context.sync();\n", "ps1": "# Commit: \"Added nuget restor & inline code\"\n#\n# Adds a new character style \"InlineCodeSnippet\" - based on the existing\n# \"CodeSnippetZchn\" (the linked character style of the paragraph style\n# \"CodeSnippet\") - so inline code spans can be formatted consistently with\n# the block CodeSnippet style. Mirrors:\n#   <w:style w:type=\"character\" w:customStyle=\"1\" w:styleId=\"InlineCodeSnippet\">\n#     <w:name w:val=\"InlineCodeSnippet\"/>\n#     <w:basedOn w:val=\"CodeSnippetZchn\"/>\n#     <w:uiPriority w:val=\"1\"/>\n#     <w:qFormat/>\n#   </w:style>\n\n$d = $word.ActiveDocument\n\n# wdStyleTypeCharacter = 2\n$inlineCodeSnippet = $d.Styles.Add(\"InlineCodeSnippet\", 2)\n\n# basedOn -> <w:basedOn w:val=\"CodeSnippetZchn\"/>\n$inlineCodeSnippet.BaseStyle = $d.Styles(\"CodeSnippetZchn\")\n# uiPriority -> <w:uiPriority w:val=\"1\"/>\n$inlineCodeSnippet.Priority = 1\n# qFormat -> <w:qFormat/>\n$inlineCodeSnippet.QuickStyle = $true\n"}
